# Append three new student records to Sheet1 (rows 4-6), matching the
# columns: studentname, fathername, mothername, gender, age, dob,
# admissioninto, address, contact_address, aadhar, contact1, contact2
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(
    @("gdgdfdg",  "dfg",    "gdg",    "male", "43",  "2024-02-02", "3rd class", "fgdffdffdsfsd", "dfsfs",  "645654654654", "5646456654", "6456546546"),
    @("Swaroop2", "dfg",    "gfdg",   "male", "654", "2024-02-13", "LKG",       "65645",         "645645", "655464564564", "6456456546", "6456456456"),
    @("jbkdsjok", "hfghgf", "hfghgf", "male", "756", "2024-02-05", "3rd class", "76756756",      "756756", "424554654645", "6465465645", "6546456546")
)

$startRow = 4
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    for ($c = 0; $c -lt $data.Count; $c++) {
        # Leading apostrophe forces the value to be stored as text (shared
        # string), matching every other cell in this sheet instead of
        # letting Excel auto-coerce numeric/date-looking strings.
        $ws.Cells.Item($r, $c + 1).Value = "'" + $data[$c]
    }
}
